$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 70
$ws.Range("C2").Value = 0.07000000000000001
$ws.Range("D2").Value = 0.0007000000000000001

$ws.Range("B3").Value = 130
$ws.Range("C3").Value = 0.07000000000000001
$ws.Range("D3").Value = 0.0007000000000000001

$ws.Range("B4").Value = 70
$ws.Range("C4").Value = 0.13
$ws.Range("D4").Value = 0.0007000000000000001

$ws.Range("B5").Value = 130
$ws.Range("C5").Value = 0.13
$ws.Range("D5").Value = 0.0007000000000000001

$ws.Range("B6").Value = 70
$ws.Range("C6").Value = 0.07000000000000001
$ws.Range("D6").Value = 0.0013

$ws.Range("B7").Value = 130
$ws.Range("C7").Value = 0.07000000000000001
$ws.Range("D7").Value = 0.0013

$ws.Range("B8").Value = 70
$ws.Range("C8").Value = 0.13
$ws.Range("D8").Value = 0.0013

$ws.Range("B9").Value = 130
$ws.Range("C9").Value = 0.13
$ws.Range("D9").Value = 0.0013
